# Apply updates described by the commit diff to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column I ("Pre_ISI") values updated to a uniform 3.01 across many rows.
$rowsForI = @(6,7,11,13,23,24,34,46,48,55,58,59,61,65,66,70,75,78,80,84,85,90,91,98,99,102,104)
foreach ($r in $rowsForI) {
    $ws.Cells.Item($r, 9).Value = 3.01
}

# Column G ("Outcome") flips on a few rows.
$ws.Cells.Item(95, 7).Value = 0
$ws.Cells.Item(101, 7).Value = 1
$ws.Cells.Item(102, 7).Value = 0

# Update the active selection to match the saved cursor position.
$ws.Range("P11").Select()
